$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add a new "Sq Class to Struct" column (H) to the legal-move stats
# table, mirroring the formatting of the previous last column (G),
# then strip the now-interior column G of its "last column" right
# border so H becomes the new right-most, bordered column.
# ------------------------------------------------------------------

# 1) Clone G's look (number formats, fills, borders, alignment) onto H
#    before touching any values - this gives H the "last column" style
#    (including the heavier right-hand border) for free.
$ws.Range("G1:G10").Copy()
$ws.Range("H1:H10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) G is no longer the last column, drop its right-hand medium border
$ws.Range("G1:G10").Borders.Item(10).LineStyle = -4142

# 3) Column width for the new column, matching the source sheet
$ws.Columns.Item(8).ColumnWidth = 18.140625

# 4) Header text
$ws.Range("H1").Value = "Sq Class to Struct"

# 5) Formula + data values for the new column
$ws.Range("H2").Formula = "=B8/H8"

$ws.Range("H3").Value = 0.001
$ws.Range("H4").Value = 0.0030000000000000001
$ws.Range("H5").Value = 0.070000000000000007
$ws.Range("H6").Value = 1.9590000000000001
$ws.Range("H7").Value = 49.692999999999998
$ws.Range("H8").Value = 1228.047
# H9 and H10 intentionally stay blank (matches G9/G10 pattern)

# 6) Selection parity with the authored workbook
$ws.Range("H3").Select()

$wb.Save()
